# Add "proximal and distal arm" parameters (v1, no gear) to the CAD parameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows in the middle of the existing table -------------
# New row 16: turntable_retention_shaft_length (pushes spine_* / shoulder_* down)
$ws.Rows.Item(16).Insert()
# New row 23: shoulder_gear_thickness (after shoulder_gear_diameter which is now row 22)
$ws.Rows.Item(23).Insert()

# --- (Re)write every data row from 16 downward so the whole block is correct ---
$rows = @(
    @(16, "turntable_retention_shaft_length", 6),
    @(17, "spine_diameter", 1.25),
    @(18, "spine_height", 1),
    @(19, "shoulder_length", 2),
    @(20, "shoulder_width", 2.75),
    @(21, "shoulder_wall_thickness", 0.1875),
    @(22, "shoulder_gear_diameter", 2.5),
    @(23, "shoulder_gear_thickness", 0.25),
    @(24, "shoulder_motor_gear_diameter", 1),
    @(25, "shoulder_gear_shaft_diameter", 0.5),
    @(26, "shoulder_motor_boss", 1),
    @(27, "arm_proximal_length", 8),
    @(28, "arm_proximal_thickness", 0.375),
    @(29, "arm_proximal_central_diameter", 1),
    @(30, "arm_proximal_motor_boss", 1),
    @(31, "arm_distal_length", 6)
)

# Write the new labels first, in the same order the author originally typed them
# (turntable shaft length, shoulder gear thickness, then the arm fields -- proximal
# length & distal length together, followed by the rest of the proximal fields).
$labelOrder = @(16, 23, 27, 31, 29, 28, 30)
foreach ($r in $labelOrder) {
    foreach ($row in $rows) {
        if ($row[0] -eq $r) {
            $ws.Cells.Item($r, 1).Value = $row[1]
        }
    }
}

# Numeric format used by the existing "value" column (fraction style, e.g. "# ?/?")
$valueNumberFormat = $ws.Cells.Item(15, 2).NumberFormat

foreach ($row in $rows) {
    $r = $row[0]
    $name = $row[1]
    $val = $row[2]

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $val
    $ws.Cells.Item($r, 2).NumberFormat = $valueNumberFormat
    $ws.Cells.Item($r, 3).Value = "in"
}

# --- Column widths (B narrower, C narrower) -------------------------------
$ws.Columns.Item(2).ColumnWidth = 5.8333333
$ws.Columns.Item(3).ColumnWidth = 3.8333333

# --- View state: selection on B30, scrolled so row 23 is at top ----------
$ws.Range("B30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1

# --- Workbook window geometry (best effort) -------------------------------
$win.Left = 4620
$win.Top = 1365
$win.Width = 15375
$win.Height = 7995
